$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The roster rows for "Anthony Lamb (TW)" (row 8) and "Moses Moody" (row 9)
# need to swap places, so Moses Moody ends up listed before Anthony Lamb
# (TW). Swap the full data (columns B:K - the "No." column A is left as-is
# since it is just the running index) between the two rows using a scratch
# area so that data types (numbers vs. text) are preserved exactly.

$srcRow = 8
$dstRow = 9
$scratchRow = 100

$rngSrc = $ws.Range("B" + $srcRow + ":K" + $srcRow)
$rngDst = $ws.Range("B" + $dstRow + ":K" + $dstRow)
$rngScratch = $ws.Range("B" + $scratchRow + ":K" + $scratchRow)

$rngSrc.Copy($rngScratch)
$rngDst.Copy($rngSrc)
$rngScratch.Copy($rngDst)
$rngScratch.Clear()
